$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.402.74'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '3.681.86'
$ws.Range("E3").Value = '  -2.42%  '
$ws.Range("E4").Value = '  +0.13%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '682.42'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.76%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '159.52'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -5.02%  '
$ws.Range("D7").Value = '3.681.13'
$ws.Range("E7").Value = '  -2.43%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  -5.04%  '
$ws.Range("E10").Value = '  -8.45%  '
$ws.Range("E11").Value = '  -5.39%  '
$ws.Range("E12").Value = '  -8.12%  '
$ws.Range("E13").Value = '  -5.96%  '
$ws.Range("D14").Value = '4.302.64'
$ws.Range("E14").Value = '  -2.33%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '32.36'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -9.15%  '
$ws.Range("D16").Value = '3.690.01'
$ws.Range("E16").Value = '  -1.65%  '
$ws.Range("D17").Value = '69.358.70'
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("E18").Value = '  -0.46%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '15.79'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -9.05%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.44'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -9.36%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '468.10'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -8.68%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '10.06'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.47%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.645'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -8.97%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '79.61'
$c.Style = "Normal"
$ws.Range("D25").Value = '3.829.97'
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -12.50%  '
$ws.Range("E28").Value = '  -12.19%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.14'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -9.57%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.68'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -7.87%  '
$ws.Range("E31").Value = '  -11.09%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.00'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -9.64%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '6.55'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -9.43%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -7.37%  '
$ws.Range("E36").Value = '  -3.80%  '
$ws.Range("D37").Value = '3.655.93'
$ws.Range("E37").Value = '  -2.08%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '8.10'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -11.66%  '
$ws.Range("E39").Value = '  -6.07%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.27'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.13%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0897'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -9.08%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("E44").Value = '  -6.08%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '164.49'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.57%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '47.55'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.85%  '
$ws.Range("E47").Value = '  -1.41%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -12.94%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.29'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -5.48%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.000274'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -7.90%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '28.09'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.60%  '
